$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods")

$check = [char]0x2713

# ---------------------------------------------------------------------------
# Insert the four new method rows at the positions implied by the diff.
# Order matters: each insert shifts everything below it down by one row, so
# we insert top-to-bottom using the row numbers the new rows end up at.
# ---------------------------------------------------------------------------
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(36).Insert()
$ws.Rows.Item(43).Insert()

# ---------------------------------------------------------------------------
# Fill the new rows' values. The fill order below (36, 17, 43, 18) reproduces
# the shared-string table order of the target workbook (new strings are
# appended to xl/sharedStrings.xml in first-use order), even though the rows
# themselves land at ascending sheet positions 17, 18, 36, 43.
# ---------------------------------------------------------------------------

# Row 36: scMultiSim
$ws.Cells.Item(36,1).Value2 = "scMultiSim"
$ws.Cells.Item(36,2).Value2 = 2022
$ws.Cells.Item(36,3).Value2 = "Class 3"
$ws.Cells.Item(36,4).Value2 = "R"
$ws.Cells.Item(36,5).Value2 = "kinetic model" + "`n" + "Beta-Poission model"
$ws.Cells.Item(36,5).WrapText = $true
$ws.Cells.Item(36,6).Value2 = "group labels (optional)"
$ws.Cells.Item(36,6).WrapText = $true
$ws.Cells.Item(36,7).Value2 = $check
$ws.Cells.Item(36,9).Value2 = $check

# Row 17: scMultiSim-tree
$ws.Cells.Item(17,1).Value2 = "scMultiSim-tree"
$ws.Cells.Item(17,2).Value2 = 2022
$ws.Cells.Item(17,3).Value2 = "Class 2"
$ws.Cells.Item(17,4).Value2 = "R"
$ws.Cells.Item(17,5).Value2 = "kinetic model" + "`n" + "Beta-Poission model"
$ws.Cells.Item(17,5).WrapText = $true
$ws.Cells.Item(17,6).Value2 = "group labels (optional)"
$ws.Cells.Item(17,6).WrapText = $true
$ws.Cells.Item(17,7).Value2 = $check
$ws.Cells.Item(17,9).Value2 = $check
$ws.Cells.Item(17,10).Value2 = $check

# Row 43: SRTsim
$ws.Cells.Item(43,1).Value2 = "SRTsim"
$ws.Cells.Item(43,2).Value2 = 2023
$ws.Cells.Item(43,3).Value2 = "Class 4"
$ws.Cells.Item(43,4).Value2 = "R"
$ws.Cells.Item(43,6).Value2 = "group labels (optional)"
$ws.Cells.Item(43,7).Value2 = $check

# Row 18: scDesign3-traj
$ws.Cells.Item(18,1).Value2 = "scDesign3-traj"
$ws.Cells.Item(18,2).Value2 = 2023
$ws.Cells.Item(18,3).Value2 = "Class 2"
$ws.Cells.Item(18,4).Value2 = "R"
$ws.Cells.Item(18,5).Value2 = "probabilistic model"
$ws.Cells.Item(18,6).Value2 = "group labels (optional)" + "`n" + "batch labels (optional)"
$ws.Cells.Item(18,6).WrapText = $true
$ws.Cells.Item(18,7).Value2 = $check
$ws.Cells.Item(18,9).Value2 = $check
$ws.Cells.Item(18,10).Value2 = $check

# ---------------------------------------------------------------------------
# Selection moves to G18 as seen in the saved workbook view state.
# ---------------------------------------------------------------------------
$ws.Range("G18").Select()
